# Refresh cryptocurrency price/volume figures (Price = column D, Volume(1h) = column E).
# Numeric-looking Price strings are entered with a leading apostrophe so Excel keeps
# them as literal text (preserving things like trailing zeros, e.g. "7.00" or "0.150")
# instead of silently re-parsing them into a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.107.75'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '3.148.41'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''604.65'
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("D6").Value = '''143.53'
$ws.Range("E6").Value = '  -3.31%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.143.13'
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").Value = '''0.524'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").Value = '''0.150'
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").Value = '''5.41'
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("D14").Value = '''35.14'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '3.658.88'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = '64.079.59'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '3.139.56'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").Value = '''6.84'
$ws.Range("E19").Value = '  -1.60%  '
$ws.Range("D20").Value = '''486.41'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").Value = '''14.73'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '''0.709'
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("D23").Value = '''7.65'
$ws.Range("E23").Value = '  -4.87%  '
$ws.Range("D24").Value = '''87.13'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''2.76'
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("D28").Value = '''8.25'
$ws.Range("E28").Value = '  -3.91%  '
$ws.Range("D29").Value = '''7.00'
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '''2.07'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").Value = '''27.19'
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("E32").Value = '  -6.56%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -3.08%  '
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("D36").Value = '''6.04'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '''52.70'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '0.0₃0746'
$ws.Range("E38").Value = '  -5.57%  '
$ws.Range("D39").Value = '''2.99'
$ws.Range("E39").Value = '  -8.09%  '
$ws.Range("D40").Value = '''438.55'
$ws.Range("E40").Value = '  -4.76%  '
$ws.Range("D41").Value = '''0.0397'
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").Value = '''0.121'
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '''8.29'
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").Value = '2.913.04'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = '''0.261'
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("D46").Value = '''2.21'
$ws.Range("E46").Value = '  -5.63%  '
$ws.Range("D47").Value = '''2.41'
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '''26.02'
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").Value = '''120.49'
$ws.Range("E51").Value = '  -0.08%  '
